$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list with GitHub Actions
# Each cell is temporarily formatted as Text ('@') before the write so
# price strings that look numeric (e.g. '22.11', '1.00') are not
# auto-coerced into floats by Excel's input parser, then ClearFormats()
# strips the temporary number format again so no new cell style lingers.

function Set-TextValue {
    param($range, [string]$value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# Column D (Price) updates
Set-TextValue $ws.Range('D2') '27.244.15'
Set-TextValue $ws.Range('D3') '1.574.74'
Set-TextValue $ws.Range('D5') '211.74'
Set-TextValue $ws.Range('D8') '22.11'
Set-TextValue $ws.Range('D9') '0.250'
Set-TextValue $ws.Range('D12') '1.796.15'
Set-TextValue $ws.Range('D13') '1.576.08'
Set-TextValue $ws.Range('D16') '27.189.01'
Set-TextValue $ws.Range('D17') '62.34'
Set-TextValue $ws.Range('D18') '7.47'
Set-TextValue $ws.Range('D20') '216.77'
Set-TextValue $ws.Range('D23') '9.26'
Set-TextValue $ws.Range('D25') '154.11'
Set-TextValue $ws.Range('D26') '6.69'
Set-TextValue $ws.Range('D27') '15.12'
Set-TextValue $ws.Range('D31') '0.0473'
Set-TextValue $ws.Range('D33') '3.18'
Set-TextValue $ws.Range('D34') '1.456.24'
Set-TextValue $ws.Range('D41') '0.811'
Set-TextValue $ws.Range('D43') '2.35'
Set-TextValue $ws.Range('D44') '1.00'
Set-TextValue $ws.Range('D45') '64.66'
Set-TextValue $ws.Range('D46') '1.74'
Set-TextValue $ws.Range('D47') '1.707.21'
Set-TextValue $ws.Range('D48') '85.87'
Set-TextValue $ws.Range('D50') '0.0525'
Set-TextValue $ws.Range('D51') '0.0962'

# Column E (Volume 1h) updates
Set-TextValue $ws.Range('E2') '  +0.92%  '
Set-TextValue $ws.Range('E3') '  +0.86%  '
Set-TextValue $ws.Range('E4') '  +0.41%  '
Set-TextValue $ws.Range('E5') '  +2.00%  '
Set-TextValue $ws.Range('E6') '  +0.71%  '
Set-TextValue $ws.Range('E7') '  +0.32%  '
Set-TextValue $ws.Range('E8') '  -0.13%  '
Set-TextValue $ws.Range('E9') '  +0.54%  '
Set-TextValue $ws.Range('E10') '  +0.73%  '
Set-TextValue $ws.Range('E11') '  +1.01%  '
Set-TextValue $ws.Range('E12') '  +0.71%  '
Set-TextValue $ws.Range('E13') '  +1.28%  '
Set-TextValue $ws.Range('E14') '  +0.84%  '
Set-TextValue $ws.Range('E15') '  +0.08%  '
Set-TextValue $ws.Range('E16') '  +0.78%  '
Set-TextValue $ws.Range('E17') '  +0.38%  '
Set-TextValue $ws.Range('E18') '  +0.91%  '
Set-TextValue $ws.Range('E19') '  -0.37%  '
Set-TextValue $ws.Range('E20') '  -0.46%  '
Set-TextValue $ws.Range('E21') '  +0.40%  '
Set-TextValue $ws.Range('E22') '  +1.23%  '
Set-TextValue $ws.Range('E24') '  +1.16%  '
Set-TextValue $ws.Range('E25') '  +0.37%  '
Set-TextValue $ws.Range('E26') '  +1.00%  '
Set-TextValue $ws.Range('E27') '  +0.35%  '
Set-TextValue $ws.Range('E28') '  +2.42%  '
Set-TextValue $ws.Range('E29') '  +0.20%  '
Set-TextValue $ws.Range('E30') '  +2.97%  '
Set-TextValue $ws.Range('E31') '  +0.96%  '
Set-TextValue $ws.Range('E32') '  +0.27%  '
Set-TextValue $ws.Range('E33') '  +1.44%  '
Set-TextValue $ws.Range('E34') '  +2.13%  '
Set-TextValue $ws.Range('E35') '  +5.46%  '
Set-TextValue $ws.Range('E36') '  +0.52%  '
Set-TextValue $ws.Range('E37') '  +1.12%  '
Set-TextValue $ws.Range('E38') '  +1.01%  '
Set-TextValue $ws.Range('E39') '  +0.56%  '
Set-TextValue $ws.Range('E40') '  +1.87%  '
Set-TextValue $ws.Range('E41') '  +0.18%  '
Set-TextValue $ws.Range('E42') '  +0.32%  '
Set-TextValue $ws.Range('E43') '  +0.93%  '
Set-TextValue $ws.Range('E44') '  -0.05%  '
Set-TextValue $ws.Range('E45') '  -0.40%  '
Set-TextValue $ws.Range('E46') '  -0.34%  '
Set-TextValue $ws.Range('E47') '  +0.62%  '
Set-TextValue $ws.Range('E48') '  -1.87%  '
Set-TextValue $ws.Range('E49') '  +3.86%  '
Set-TextValue $ws.Range('E50') '  +0.48%  '
